$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the date / price data between row 3 and row 4
# Row 3 becomes what row 4 used to hold, and vice versa.

$ws.Range("D3").Value = 44525
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 8000
$ws.Range("P3").Value = 533

$ws.Range("D4").Value = 44508
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("P4").Value = 667
